# Apply edits described by the diff:
# - Update the actual time length to complete for "Discussion question 2" (C12)
# - Total (C20) recalculates automatically via its SUM formula
# - Update active selection to C13
# - Update workbook window position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# Update C12: actual time length to complete for DQ2 row
$ws.Range("C12").Value = 0.09305555555555556

# Move the active cell selection to C13 on the week1 sheet
$ws.Activate()
$ws.Range("C13").Select()

# Update workbook window position/size
$excel.ActiveWindow.Left = 8280
$excel.ActiveWindow.Top = 5600
